$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain decimal number (e.g. "5.90").
# Excel auto-converts such text typed into Value to a real number,
# which would also attach a text-number-format style to the cell.
# Force text entry by temporarily setting the cell format to Text ("@"),
# then clear the formatting again afterwards so the cell ends up with
# no style index at all (matching the original file).
$riskyCells = @("D5", "D6", "D10", "D11", "D13", "D17", "D20", "D21", "D22", "D23", "D28", "D29", "D30", "D31", "D32", "D39", "D41", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($c in $riskyCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply the updated price (D) and volume (E) values row by row.
$ws.Range("D2").Value = '44.649.47'
$ws.Range("E2").Value = '  +0.97%  '
$ws.Range("D3").Value = '2.247.89'
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = '306.92'
$ws.Range("E5").Value = '  +0.16%  '
$ws.Range("D6").Value = '95.01'
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("E7").Value = '  -0.38%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("E9").Value = '  -1.87%  '
$ws.Range("D10").Value = '34.97'
$ws.Range("E10").Value = '  +0.52%  '
$ws.Range("D11").Value = '0.0802'
$ws.Range("E11").Value = '  -1.23%  '
$ws.Range("E12").Value = '  +0.13%  '
$ws.Range("D13").Value = '0.105'
$ws.Range("E13").Value = '  +0.18%  '
$ws.Range("D14").Value = '2.591.37'
$ws.Range("E14").Value = '  +0.20%  '
$ws.Range("D15").Value = '2.245.13'
$ws.Range("E15").Value = '  -3.58%  '
$ws.Range("E16").Value = '  +0.08%  '
$ws.Range("D17").Value = '13.56'
$ws.Range("E17").Value = '  +0.08%  '
$ws.Range("D18").Value = '44.423.81'
$ws.Range("E18").Value = '  +1.05%  '
$ws.Range("E19").Value = '  -2.93%  '
$ws.Range("D20").Value = '6.19'
$ws.Range("E20").Value = '  -3.16%  '
$ws.Range("D21").Value = '11.72'
$ws.Range("E21").Value = '  -3.34%  '
$ws.Range("D22").Value = '65.35'
$ws.Range("E22").Value = '  -0.40%  '
$ws.Range("D23").Value = '237.68'
$ws.Range("E23").Value = '  -0.57%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("E25").Value = '  -1.41%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D28").Value = '9.77'
$ws.Range("E28").Value = '  -1.71%  '
$ws.Range("D29").Value = '36.98'
$ws.Range("E29").Value = '  -3.53%  '
$ws.Range("D30").Value = '5.90'
$ws.Range("D31").Value = '20.00'
$ws.Range("E31").Value = '  -0.32%  '
$ws.Range("D32").Value = '147.89'
$ws.Range("E32").Value = '  -3.56%  '
$ws.Range("E33").Value = '  -1.49%  '
$ws.Range("E34").Value = '  +0.11%  '
$ws.Range("E35").Value = '  +0.24%  '
$ws.Range("E36").Value = '  +1.06%  '
$ws.Range("E37").Value = '  -1.51%  '
$ws.Range("E38").Value = '  +4.88%  '
$ws.Range("D39").Value = '15.22'
$ws.Range("E39").Value = '  +5.77%  '
$ws.Range("E40").Value = '  -6.17%  '
$ws.Range("D41").Value = '3.78'
$ws.Range("E41").Value = '  -1.23%  '
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("E43").Value = '  +0.09%  '
$ws.Range("D44").Value = '1.809.58'
$ws.Range("E44").Value = '  +3.77%  '
$ws.Range("E45").Value = '  +11.75%  '
$ws.Range("D46").Value = '81.86'
$ws.Range("E46").Value = '  -1.16%  '
$ws.Range("D47").Value = '0.188'
$ws.Range("E47").Value = '  -1.79%  '
$ws.Range("D48").Value = '98.43'
$ws.Range("E48").Value = '  -1.54%  '
$ws.Range("D49").Value = '4.83'
$ws.Range("E49").Value = '  -2.05%  '
$ws.Range("D50").Value = '68.97'
$ws.Range("E50").Value = '  +2.61%  '
$ws.Range("D51").Value = '54.05'
$ws.Range("E51").Value = '  -0.96%  '

# Remove the temporary text formatting so the cells keep their original
# (unset) style, now that the literal text values are locked in.
foreach ($c in $riskyCells) {
    $ws.Range($c).ClearFormats()
}
